# Fruta / hortaliza, semanal
# Insert a new data row at row 424 (pushing existing rows 424:519 down to 425:520)
# and populate it with the new weekly record. Columns A,B,C,E,F,G,H,I,J,R carry
# the same constant values used throughout this block of rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(424).Insert()

$ws.Cells.Item(424, 1).Value2 = 10
$ws.Cells.Item(424, 2).Value2 = "Vega Modelo de Temuco"
$ws.Cells.Item(424, 3).Value2 = "La Araucanía"
$ws.Cells.Item(424, 4).Value2 = 44782
$ws.Cells.Item(424, 5).Value2 = 9
$ws.Cells.Item(424, 6).Value2 = "Fruta"
$ws.Cells.Item(424, 7).Value2 = 100108
$ws.Cells.Item(424, 8).Value2 = "Tropicales y subtropicales"
$ws.Cells.Item(424, 9).Value2 = 100108005
$ws.Cells.Item(424, 10).Value2 = "Piña"
$ws.Cells.Item(424, 11).Value2 = "Caramelo"
$ws.Cells.Item(424, 12).Value2 = "Segunda"
$ws.Cells.Item(424, 13).Value2 = 65
$ws.Cells.Item(424, 14).Value2 = 22000
$ws.Cells.Item(424, 15).Value2 = 22000
$ws.Cells.Item(424, 16).Value2 = 22000
$ws.Cells.Item(424, 17).Value2 = "$/caja 14 unidades"
$ws.Cells.Item(424, 18).Value2 = "Ecuador"
$ws.Cells.Item(424, 19).Value2 = 1571
$ws.Cells.Item(424, 20).Value2 = 14
